{"js": "// Replace the stale EPackageImpl hashcode with the updated one produced\n// by the regenerated test output (both occurrences in the document body).\nconst oldText = \"org.eclipse.emf.ecore.impl.EPackageImpl@231e1f6e\";\nconst newText = \"org.eclipse.emf.ecore.impl.EPackageImpl@18a6141a\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the stale EPackageImpl hashcode with the updated one produced\n# by the regenerated test output (both occurrences in the document body).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"org.eclipse.emf.ecore.impl.EPackageImpl@231e1f6e\"\n$find.Replacement.Text = \"org.eclipse.emf.ecore.impl.EPackageImpl@18a6141a\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 1\n\n# wdFindContinue=1, wdReplaceAll=2\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
